# Product and Order Inventory Test
# Applies the order-import.xlsx template update:
#  - removes the old blank "row 5" placeholder row (Sheet2 formulas re-point accordingly)
#  - refreshes sample data in the remaining rows
#  - adds new "SKU" / "Order" sample columns (AC:AD)
#  - fixes up the hyperlink that tracked the now-shifted data row
#  - restores the plain (non-numeric-formatted) style on the weight-length cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- remove the stale placeholder row (old row 5); rows below shift up ---
$ws.Rows.Item(5).Delete()

# --- row 2 (homedelivery / fedex order) ---
$ws.Range("D2").Value = 2455
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = 12
$ws.Range("R2").Value = "Santos"
$ws.Range("S2").Value = "SP"
$ws.Range("T2").Value = "BR"
$ws.Range("W2").Value = 5
$ws.Range("Y2").Value = "Racing Car"
$ws.Range("AC2").Value = "SKU0024"
$ws.Range("AD2").Value = 30203

# --- row 3 ---
$ws.Range("D3").Value = 2455
$ws.Range("W3").Value = 5
$ws.Range("Y3").Value = "A5 Cell"
$ws.Range("AC3").Value = "SKU0023"
$ws.Range("AD3").Value = 30202

# --- row 4 ---
$ws.Range("D4").Value = 2455
$ws.Range("W4").Value = 5
$ws.Range("Y4").Value = "Charger"
$ws.Range("AA4").Value = "yes"
$ws.Range("AC4").Value = "SKU0022"
$ws.Range("AD4").Value = 30201

# --- row 5 (was row 6 pre-delete; now holds the "usps" sample order) ---
$ws.Range("A5").Value = "homedelivery"
$ws.Range("B5").Value = "usps"
$ws.Range("C5").Value = 30065626894
$ws.Range("D5").Value = 2456
$ws.Range("R5").Value = "Male"
$ws.Range("Y5").Value = "Racing Car"
$ws.Range("AB5").Value = "yes"
$ws.Range("AC5").Value = "SKU0024"
$ws.Range("AD5").Value = 30203

# --- new sample columns header ---
$ws.Range("AC1").Value = "SKU"
$ws.Range("AD1").Value = "Order"

# --- hyperlinks: drop the stale set and re-add at the correct (shifted) cells ---
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("L2"), "mailto:saqib@874@gmail.com")
$ws.Hyperlinks.Add($ws.Range("L5"), "mailto:82huasd@mail.com")

$ws.Range("A5").Select()
